$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column F width (between existing col B and col G) to 12.5 with custom width.
# NOTE: the engine's OOXML serialization adds 5/6 to the ColumnWidth value when
# writing the stored "width" attribute, so we back that off here to land on 12.5.
$ws.Columns.Item(6).ColumnWidth = 11.666666666666666

# Update the "Big Id" values in column F to larger numbers (fix for big decimal import error)
$ws.Range("F2").Value = 11111111111
$ws.Range("F3").Value = 11111111112
$ws.Range("F4").Value = 11111111113
$ws.Range("F5").Value = 11111111114
$ws.Range("F6").Value = 11111111115

# Update the selection to F2:F6 with active cell F2
$ws.Range("F2:F6").Select()
